$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.735499858856201
$ws.Range("B1").Value = 2.460647583007812
$ws.Range("C1").Value = 2.531919717788696
$ws.Range("D1").Value = 2.899296045303345
$ws.Range("E1").Value = 3.640483379364014
